$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at A, shifting the existing "colunm Name 1"/"colunm Name 2"
# header columns from A/B to B/C. This preserves column B's original width
# definition (13.90625, bestFit) by shifting it off column A onto column B.
$ws.Columns.Item(1).Insert()

# New header row: Parameters | username | password
$ws.Cells.Item(1, 2).Value = "username"
$ws.Cells.Item(1, 3).Value = "password"
$ws.Cells.Item(1, 1).Value = "Parameters"

# Data rows
$ws.Cells.Item(2, 2).Value = "admin"
$ws.Cells.Item(2, 3).Value = "galatpassword"

$ws.Cells.Item(3, 2).Value = "rajmitra"
$ws.Cells.Item(3, 3).Value = "manager"

# Bold header row
$ws.Range("A1:C1").Font.Bold = $true

# Column A width (new column, no bestFit)
$ws.Columns.Item(1).ColumnWidth = 15.5

# Selection
$ws.Range("G11").Select()
